# gra2005 + gra2009 new storageType
# Inserts a new "doNotDisplayThousandsSeparator" column right after "name"
# (i.e. before the old column B / "panelIdentifier") and fills it in:
# begin_y / end_y (gra2005 / gra2009 rows) get "true", everything else "false".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing columns B:N one column to the right, to make room for
# the new "doNotDisplayThousandsSeparator" column.
$ws.Columns("B").Insert()

# New header
$ws.Range("B1").Value = "doNotDisplayThousandsSeparator"

# Default all data rows (2-8) to "false" ...
$ws.Range("B2:B8").Value = "false"

# ... except begin_y (row 6) and end_y (row 8), which are "true".
$ws.Range("B6").Value = "true"
$ws.Range("B8").Value = "true"

# Match the saved selection state recorded in the workbook.
$ws.Range("B7").Select()
